$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$req1 = "LOT2052 -  Tecnologia de Bebidas Experimental  (Indicação de Conjunto)`n"
$req2 = "LOT2028 -  Tecnologia de Processos Fermentativos  (Requisito fraco)`n"

$ws.Range("B24").Value = $req2
$ws.Range("C24").Value = $req2
$ws.Range("B25").Value = $req1
$ws.Range("C25").Value = $req1
